# "tilføjet alle opgs til samlet dokument"
# Adds all remaining "mangler" (missing) task notes to the consolidated
# overview sheet, plus a small red/orange/green status legend, and removes
# the now-obsolete "en masse spg" note under L06 - ANN.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# BGR packed colors (as used by the Excel COM object model):
#   green  FF00B050 -> done
#   orange FFFFC000 -> partially missing / needs work
#   red    FFFF0000 -> missing entirely
$colorGreen  = 5287936   # 0x0050B0 -> RGB(0,176,80)
$colorOrange = 49407     # 0x00C0FF -> RGB(255,192,0)
$colorRed    = 255       # 0x0000FF -> RGB(255,0,0)

# 1) Drop the stray "en masse spg" note under L06 - ANN (row 18) -- it is
#    superseded by the fuller breakdown added below.
$ws.Range("C18").ClearContents()

# 2) "Mangler alt" under L05 - Train linear regression is still fully missing
#    -> flag it red.
$ws.Range("C17").Interior.Color = $colorRed

# 3) Fill out the L06 - ANN section with the per-question status notes.
$ws.Range("D20").Value = "mangler intro"
$ws.Range("D20").Interior.Color = $colorOrange

$ws.Range("C21").Value = "mangler tekst"
$ws.Range("C21").Interior.Color = $colorOrange

# 4) Small status legend next to the top section (L03).
$ws.Range("E3").Value = "mangler kode/opgave"
$ws.Range("E3").Interior.Color = $colorRed

$ws.Range("E4").Value = "mangler forklaring"
$ws.Range("E4").Interior.Color = $colorOrange

$ws.Range("E5").Value = "done"
$ws.Range("E5").Interior.Color = $colorGreen

# 5) Remaining L06 - ANN question notes.
$ws.Range("C22").Value = "husk at indsæt nyt billede med endelig data THOMAS"
$ws.Range("C22").Interior.Color = $colorOrange

$ws.Range("C23").Value = "mangler tekst, ryk kode op I markdown spørg lasse"
$ws.Range("C23").Interior.Color = $colorOrange

$ws.Range("C24").Value = "graf er skæv, lav om???? Passer ikke"
$ws.Range("C24").Interior.Color = $colorRed

$ws.Range("C25").Value = "mangler kode, victoria har det"
$ws.Range("C25").Interior.Color = $colorOrange

$ws.Range("C26").Value = "mangler tekst"
$ws.Range("C26").Interior.Color = $colorOrange

# 6) Widen the columns to fit the new text.
# (Input values below are chosen so the engine's pixel-quantized
# ColumnWidth conversion lands as close as possible to the authored
# widths of 42.88671875 / 21.21875 / 19.33203125 characters.)
$ws.Columns("C").ColumnWidth = 41.92
$ws.Columns("D").ColumnWidth = 20.25
$ws.Columns("E").ColumnWidth = 18.42

# 7) Leave the selection where the user last clicked while editing.
$ws.Range("F21").Select() | Out-Null
